# Auto-generated edit script: updates cryptos list figures (price/volume)
# and fixes the RenderToken/VeChain row ordering, per commit
# "Updated cryptos list on Fri Jun 16 05:08:25 UTC 2023 with GitHub Actions".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.496.60'
$ws.Range("E2").Value = '  +1.86%  '
$ws.Range("D3").Value = '1.667.65'
$ws.Range("E3").Value = '  +1.05%  '
$ws.Range("D4").Value = '''0.9992'
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = '''237.56'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.40%  '
$ws.Range("E6").Value = '  +0.00%  '
$ws.Range("D7").Value = '''0.4796'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.82%  '
$ws.Range("E8").Value = '  +0.48%  '
$ws.Range("D9").Value = '''0.06171'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.57%  '
$ws.Range("D10").Value = '''0.07068'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.73%  '
$ws.Range("D11").Value = '1.667.25'
$ws.Range("E11").Value = '  +1.08%  '
$ws.Range("D12").Value = '''14.81'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.08%  '
$ws.Range("D13").Value = '''0.5882'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -5.50%  '
$ws.Range("E14").Value = '  -4.26%  '
$ws.Range("D15").Value = '''74.91'
$ws.Range("D15").Style = "Normal"
$ws.Range("E16").Value = '  +0.02%  '
$ws.Range("D17").Value = '''1.000'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.07%  '
$ws.Range("D18").Value = '25.496.00'
$ws.Range("E18").Value = '  +1.93%  '
$ws.Range("D19").Value = '''0.000006749'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.83%  '
$ws.Range("E20").Value = '  -0.77%  '
$ws.Range("D21").Value = '1.879.07'
$ws.Range("E21").Value = '  +1.00%  '
$ws.Range("D22").Value = '''4.433'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.84%  '
$ws.Range("D23").Value = '''8.729'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.19%  '
$ws.Range("E24").Value = '  -0.69%  '
$ws.Range("D25").Value = '''136.75'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.43%  '
$ws.Range("D26").Value = '''15.01'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.67%  '
$ws.Range("D27").Value = '''1.388'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.76%  '
$ws.Range("D28").Value = '''104.91'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.81%  '
$ws.Range("D29").Value = '''1.715'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +2.15%  '
$ws.Range("D30").Value = '''3.952'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +4.93%  '
$ws.Range("D31").Value = '''0.07785'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.36%  '
$ws.Range("D32").Value = '''3.643'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.45%  '
$ws.Range("D33").Value = '''0.9994'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.01%  '
$ws.Range("D34").Value = '''0.04207'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -7.44%  '
$ws.Range("E35").Value = '  +0.37%  '
$ws.Range("D36").Value = '''0.6093'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.06%  '
$ws.Range("D37").Value = '''0.9477'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.29%  '
$ws.Range("D38").Value = '''2.593'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.19%  '
$ws.Range("D39").Value = '''0.8606'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.98%  '
$ws.Range("E40").Value = '  +0.14%  '
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").Value = '''0.01472'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.89%  '
$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D42").Value = '''1.843'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.12%  '
$ws.Range("D43").Value = '''97.18'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.35%  '
$ws.Range("D44").Value = '''0.3758'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.49%  '
$ws.Range("D45").Value = '''4.834'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.87%  '
$ws.Range("D46").Value = '''0.1118'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.80%  '
$ws.Range("D47").Value = '''6.196'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.88%  '
$ws.Range("D48").Value = '''0.05256'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.10%  '
$ws.Range("D49").Value = '''29.70'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.61%  '
$ws.Range("E50").Value = '  +0.12%  '
$ws.Range("D51").Value = '''1.000'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.03%  '
